$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ah = New-Object "object[,]" 41,8
$ah[0,0] = "poorly"
$ah[0,1] = 1
$ah[0,2] = 46
$ah[0,3] = 46
$ah[0,4] = 0
$ah[0,5] = 1
$ah[0,6] = $false
$ah[0,7] = 0
$ah[1,0] = "disappointing"
$ah[1,1] = 0.7954545454545454
$ah[1,2] = 35
$ah[1,3] = 35
$ah[1,4] = 0
$ah[1,5] = 1
$ah[1,6] = $false
$ah[1,7] = 9
$ah[2,0] = "poor"
$ah[2,1] = 0.7746478873239436
$ah[2,2] = 55
$ah[2,3] = 55
$ah[2,4] = 0
$ah[2,5] = 1
$ah[2,6] = $false
$ah[2,7] = 16
$ah[3,0] = "disappointed"
$ah[3,1] = 0.7365591397849462
$ah[3,2] = 137
$ah[3,3] = 137
$ah[3,4] = 0
$ah[3,5] = 1
$ah[3,6] = $false
$ah[3,7] = 49
$ah[4,0] = "however"
$ah[4,1] = 0.734375
$ah[4,2] = 47
$ah[4,3] = 47
$ah[4,4] = 0
$ah[4,5] = 1
$ah[4,6] = $false
$ah[4,7] = 17
$ah[5,0] = "broke"
$ah[5,1] = 0.7281553398058253
$ah[5,2] = 150
$ah[5,3] = 150
$ah[5,4] = 0
$ah[5,5] = 1
$ah[5,6] = $false
$ah[5,7] = 56
$ah[6,0] = "waste"
$ah[6,1] = 0.6959459459459459
$ah[6,2] = 103
$ah[6,3] = 103
$ah[6,4] = 0
$ah[6,5] = 1
$ah[6,6] = $false
$ah[6,7] = 45
$ah[7,0] = "instead"
$ah[7,1] = 0.6875
$ah[7,2] = 33
$ah[7,3] = 33
$ah[7,4] = 0
$ah[7,5] = 1
$ah[7,6] = $false
$ah[7,7] = 15
$ah[8,0] = "junk"
$ah[8,1] = 0.6363636363636364
$ah[8,2] = 35
$ah[8,3] = 35
$ah[8,4] = 0
$ah[8,5] = 1
$ah[8,6] = $false
$ah[8,7] = 20
$ah[9,0] = "smaller"
$ah[9,1] = 0.5882352941176471
$ah[9,2] = 70
$ah[9,3] = 70
$ah[9,4] = 0
$ah[9,5] = 1
$ah[9,6] = $false
$ah[9,7] = 49
$ah[10,0] = "small"
$ah[10,1] = 0.518840579710145
$ah[10,2] = 179
$ah[10,3] = 179
$ah[10,4] = 0
$ah[10,5] = 1
$ah[10,6] = $false
$ah[10,7] = 166
$ah[11,0] = "broken"
$ah[11,1] = 0.5180722891566265
$ah[11,2] = 43
$ah[11,3] = 43
$ah[11,4] = 0
$ah[11,5] = 1
$ah[11,6] = $false
$ah[11,7] = 40
$ah[12,0] = "paint"
$ah[12,1] = 0.4761904761904762
$ah[12,2] = 30
$ah[12,3] = 30
$ah[12,4] = 0
$ah[12,5] = 1
$ah[12,6] = $false
$ah[12,7] = 33
$ah[13,0] = "di"
$ah[13,1] = 0.453125
$ah[13,2] = 29
$ah[13,3] = 29
$ah[13,4] = 0
$ah[13,5] = 1
$ah[13,6] = $false
$ah[13,7] = 35
$ah[14,0] = "apart"
$ah[14,1] = 0.4421052631578947
$ah[14,2] = 42
$ah[14,3] = 42
$ah[14,4] = 0
$ah[14,5] = 1
$ah[14,6] = $false
$ah[14,7] = 53
$ah[15,0] = "plastic"
$ah[15,1] = 0.4094488188976378
$ah[15,2] = 52
$ah[15,3] = 52
$ah[15,4] = 0
$ah[15,5] = 1
$ah[15,6] = $false
$ah[15,7] = 75
$ah[16,0] = "ok"
$ah[16,1] = 0.3984375
$ah[16,2] = 51
$ah[16,3] = 51
$ah[16,4] = 0
$ah[16,5] = 1
$ah[16,6] = $false
$ah[16,7] = 77
$ah[17,0] = "difficult"
$ah[17,1] = 0.3820224719101123
$ah[17,2] = 34
$ah[17,3] = 34
$ah[17,4] = 0
$ah[17,5] = 1
$ah[17,6] = $false
$ah[17,7] = 55
$ah[18,0] = "thought"
$ah[18,1] = 0.3415841584158416
$ah[18,2] = 69
$ah[18,3] = 69
$ah[18,4] = 0
$ah[18,5] = 1
$ah[18,6] = $false
$ah[18,7] = 133
$ah[19,0] = "bit"
$ah[19,1] = 0.3265306122448979
$ah[19,2] = 32
$ah[19,3] = 32
$ah[19,4] = 0
$ah[19,5] = 1
$ah[19,6] = $false
$ah[19,7] = 66
$ah[20,0] = "cheap"
$ah[20,1] = 0.2748815165876777
$ah[20,2] = 58
$ah[20,3] = 58
$ah[20,4] = 0
$ah[20,5] = 1
$ah[20,6] = $false
$ah[20,7] = 153
$ah[21,0] = "size"
$ah[21,1] = 0.2628865979381443
$ah[21,2] = 51
$ah[21,3] = 51
$ah[21,4] = 0
$ah[21,5] = 1
$ah[21,6] = $false
$ah[21,7] = 143
$ah[22,0] = "though"
$ah[22,1] = 0.2564102564102564
$ah[22,2] = 30
$ah[22,3] = 30
$ah[22,4] = 0
$ah[22,5] = 1
$ah[22,6] = $false
$ah[22,7] = 87
$ah[23,0] = "back"
$ah[23,1] = 0.2428571428571429
$ah[23,2] = 34
$ah[23,3] = 34
$ah[23,4] = 0
$ah[23,5] = 1
$ah[23,6] = $false
$ah[23,7] = 106
$ah[24,0] = "would"
$ah[24,1] = 0.2303120356612184
$ah[24,2] = 155
$ah[24,3] = 156
$ah[24,4] = 0.01
$ah[24,5] = 0.99
$ah[24,6] = $true
$ah[24,7] = 518
$ah[25,0] = "work"
$ah[25,1] = 0.2056962025316456
$ah[25,2] = 65
$ah[25,3] = 65
$ah[25,4] = 0
$ah[25,5] = 1
$ah[25,6] = $false
$ah[25,7] = 251
$ah[26,0] = "money"
$ah[26,1] = 0.2025316455696203
$ah[26,2] = 64
$ah[26,3] = 64
$ah[26,4] = 0
$ah[26,5] = 1
$ah[26,6] = $false
$ah[26,7] = 252
$ah[27,0] = "item"
$ah[27,1] = 0.1920289855072464
$ah[27,2] = 53
$ah[27,3] = 53
$ah[27,4] = 0
$ah[27,5] = 1
$ah[27,6] = $false
$ah[27,7] = 223
$ah[28,0] = "could"
$ah[28,1] = 0.1847133757961783
$ah[28,2] = 29
$ah[28,3] = 29
$ah[28,4] = 0
$ah[28,5] = 1
$ah[28,6] = $false
$ah[28,7] = 128
$ah[29,0] = "product"
$ah[29,1] = 0.1828193832599119
$ah[29,2] = 83
$ah[29,3] = 83
$ah[29,4] = 0
$ah[29,5] = 1
$ah[29,6] = $false
$ah[29,7] = 371
$ah[30,0] = "hard"
$ah[30,1] = 0.18
$ah[30,2] = 36
$ah[30,3] = 36
$ah[30,4] = 0
$ah[30,5] = 1
$ah[30,6] = $false
$ah[30,7] = 164
$ah[31,0] = "used"
$ah[31,1] = 0.1714285714285714
$ah[31,2] = 30
$ah[31,3] = 30
$ah[31,4] = 0
$ah[31,5] = 1
$ah[31,6] = $false
$ah[31,7] = 145
$ah[32,0] = "price"
$ah[32,1] = 0.170028818443804
$ah[32,2] = 59
$ah[32,3] = 60
$ah[32,4] = 0.02
$ah[32,5] = 0.98
$ah[32,6] = $true
$ah[32,7] = 288
$ah[33,0] = "better"
$ah[33,1] = 0.1448598130841121
$ah[33,2] = 31
$ah[33,3] = 31
$ah[33,4] = 0
$ah[33,5] = 1
$ah[33,6] = $false
$ah[33,7] = 183
$ah[34,0] = "3"
$ah[34,1] = 0.1341463414634146
$ah[34,2] = 33
$ah[34,3] = 35
$ah[34,4] = 0.06
$ah[34,5] = 0.94
$ah[34,6] = $true
$ah[34,7] = 213
$ah[35,0] = "use"
$ah[35,1] = 0.1315068493150685
$ah[35,2] = 48
$ah[35,3] = 48
$ah[35,4] = 0
$ah[35,5] = 1
$ah[35,6] = $false
$ah[35,7] = 317
$ah[36,0] = "2"
$ah[36,1] = 0.1161048689138577
$ah[36,2] = 31
$ah[36,3] = 31
$ah[36,4] = 0
$ah[36,5] = 1
$ah[36,6] = $false
$ah[36,7] = 236
$ah[37,0] = "like"
$ah[37,1] = 0.09225700164744646
$ah[37,2] = 56
$ah[37,3] = 57
$ah[37,4] = 0.02
$ah[37,5] = 0.98
$ah[37,6] = $true
$ah[37,7] = 551
$ah[38,0] = "little"
$ah[38,1] = 0.08908685968819599
$ah[38,2] = 40
$ah[38,3] = 40
$ah[38,4] = 0
$ah[38,5] = 1
$ah[38,6] = $false
$ah[38,7] = 409
$ah[39,0] = "buy"
$ah[39,1] = 0.08732394366197183
$ah[39,2] = 31
$ah[39,3] = 31
$ah[39,4] = 0
$ah[39,5] = 1
$ah[39,6] = $false
$ah[39,7] = 324
$ah[40,0] = "one"
$ah[40,1] = 0.04683544303797468
$ah[40,2] = 37
$ah[40,3] = 41
$ah[40,4] = 0.1
$ah[40,5] = 0.9
$ah[40,6] = $true
$ah[40,7] = 753
$ws.Range("A3:H43").Value = $ah

$jq = New-Object "object[,]" 15,8
$jq[0,0] = "wonderful"
$jq[0,1] = 0.8571428571428571
$jq[0,2] = 48
$jq[0,3] = 48
$jq[0,4] = 1
$jq[0,5] = 0
$jq[0,6] = $false
$jq[0,7] = 8
$jq[1,0] = "awesome"
$jq[1,1] = 0.8307692307692308
$jq[1,2] = 54
$jq[1,3] = 54
$jq[1,4] = 1
$jq[1,5] = 0
$jq[1,6] = $false
$jq[1,7] = 11
$jq[2,0] = "favorite"
$jq[2,1] = 0.6344086021505376
$jq[2,2] = 59
$jq[2,3] = 59
$jq[2,4] = 1
$jq[2,5] = 0
$jq[2,6] = $false
$jq[2,7] = 34
$jq[3,0] = "excellent"
$jq[3,1] = 0.5625
$jq[3,2] = 36
$jq[3,3] = 36
$jq[3,4] = 1
$jq[3,5] = 0
$jq[3,6] = $false
$jq[3,7] = 28
$jq[4,0] = "classic"
$jq[4,1] = 0.5471698113207547
$jq[4,2] = 29
$jq[4,3] = 29
$jq[4,4] = 1
$jq[4,5] = 0
$jq[4,6] = $false
$jq[4,7] = 24
$jq[5,0] = "great"
$jq[5,1] = 0.3713114754098361
$jq[5,2] = 453
$jq[5,3] = 453
$jq[5,4] = 1
$jq[5,5] = 0
$jq[5,6] = $false
$jq[5,7] = 767
$jq[6,0] = "love"
$jq[6,1] = 0.321377331420373
$jq[6,2] = 224
$jq[6,3] = 224
$jq[6,4] = 1
$jq[6,5] = 0
$jq[6,6] = $false
$jq[6,7] = 473
$jq[7,0] = "loves"
$jq[7,1] = 0.3008298755186722
$jq[7,2] = 145
$jq[7,3] = 145
$jq[7,4] = 1
$jq[7,5] = 0
$jq[7,6] = $false
$jq[7,7] = 337
$jq[8,0] = "best"
$jq[8,1] = 0.25
$jq[8,2] = 30
$jq[8,3] = 30
$jq[8,4] = 1
$jq[8,5] = 0
$jq[8,6] = $false
$jq[8,7] = 90
$jq[9,0] = "perfect"
$jq[9,1] = 0.2409638554216867
$jq[9,2] = 40
$jq[9,3] = 40
$jq[9,4] = 1
$jq[9,5] = 0
$jq[9,6] = $false
$jq[9,7] = 126
$jq[10,0] = "loved"
$jq[10,1] = 0.2171253822629969
$jq[10,2] = 71
$jq[10,3] = 71
$jq[10,4] = 1
$jq[10,5] = 0
$jq[10,6] = $false
$jq[10,7] = 256
$jq[11,0] = "friends"
$jq[11,1] = 0.1904761904761905
$jq[11,2] = 36
$jq[11,3] = 36
$jq[11,4] = 1
$jq[11,5] = 0
$jq[11,6] = $false
$jq[11,7] = 153
$jq[12,0] = "christmas"
$jq[12,1] = 0.1164658634538153
$jq[12,2] = 29
$jq[12,3] = 29
$jq[12,4] = 1
$jq[12,5] = 0
$jq[12,6] = $false
$jq[12,7] = 220
$jq[13,0] = "fun"
$jq[13,1] = 0.1052631578947368
$jq[13,2] = 120
$jq[13,3] = 121
$jq[13,4] = 0.99
$jq[13,5] = 0.01000000000000001
$jq[13,6] = $true
$jq[13,7] = 1020
$jq[14,0] = "game"
$jq[14,1] = 0.06627680311890838
$jq[14,2] = 102
$jq[14,3] = 104
$jq[14,4] = 0.98
$jq[14,5] = 0.02000000000000002
$jq[14,6] = $true
$jq[14,7] = 1437
$ws.Range("J3:Q17").Value = $jq
